$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: row -> (new DAMSLTag, new DialogAct)
$updates = @{
    5   = @("ba", "Appreciation")
    8   = @("b",  "Acknowledge (Backchannel)")
    20  = @("sd", "Statement-non-opinion")
    22  = @("b",  "Acknowledge (Backchannel)")
    28  = @("b",  "Acknowledge (Backchannel)")
    51  = @("sd", "Statement-non-opinion")
    55  = @("sd", "Statement-non-opinion")
    56  = @("b",  "Acknowledge (Backchannel)")
    58  = @("ba", "Appreciation")
    87  = @("b",  "Acknowledge (Backchannel)")
    89  = @("b",  "Acknowledge (Backchannel)")
    93  = @("%",  "Uninterpretable")
    94  = @("b",  "Acknowledge (Backchannel)")
    100 = @("sd", "Statement-non-opinion")
    103 = @("aa", "Agree/Accept")
    108 = @("ba", "Appreciation")
    113 = @("sd", "Statement-non-opinion")
}

foreach ($row in $updates.Keys) {
    $vals = $updates[$row]
    $ws.Range("I$row").Value = $vals[0]
    $ws.Range("J$row").Value = $vals[1]
}
